# testdata.xlsx: flip the "execute" flag for the newTest cases from "no" to
# "yes" (RUNMANAGER!C2 and DATA!B4:B5), which drops "no" from the shared
# string table entirely since nothing references it any more.

$wb = $excel.ActiveWorkbook

$wsRun  = $wb.Worksheets.Item("RUNMANAGER")
$wsData = $wb.Worksheets.Item("DATA")

# RUNMANAGER: loginLogoutTest row's "execute" column -> yes
$wsRun.Range("C2").Value = "yes"

# DATA: both newTest rows' "execute" column -> yes
$wsData.Range("B4").Value = "yes"
$wsData.Range("B5").Value = "yes"

# Leave the selection on DATA where the work happened, then switch back to
# RUNMANAGER as the active/visible tab.
$wsData.Activate()
$wsData.Range("F6").Select()

$wsRun.Activate()
$wsRun.Range("F4").Select()
